# Update the "Handoff/Handback" timestamp cells on the Overview, zh-cn and
# de-de sheets to reflect a newly (re-)generated handback report.
#
# Overview!G3  (shared with de-de!H3): 2016-09-06 00:50:11 -> 2016-09-06 00:51:08
# zh-cn!H3:                            2016-09-06 00:49:58 -> 2016-09-06 00:50:58
# zh-cn!K3:                            2016-09-06 00:50:31 -> 2016-09-06 00:51:30
# de-de!K3:                            2016-09-06 00:50:39 -> 2016-09-06 00:51:38

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet, row for 0df4348b-5154-4cbc-915c-7b38bbdec900.md
$wsOverview.Range("G3").Value = "2016-09-06 00:51:08"

# zh-cn sheet, row for 0df4348b-5154-4cbc-915c-7b38bbdec900...zh-cn.xlf
$wsZhCn.Range("H3").Value = "2016-09-06 00:50:58"
$wsZhCn.Range("K3").Value = "2016-09-06 00:51:30"

# de-de sheet, row for 0df4348b-5154-4cbc-915c-7b38bbdec900...de-de.xlf
$wsDeDe.Range("K3").Value = "2016-09-06 00:51:38"
